# Organiza los splines y las celdas de la seccion 3
# The sheet was a 3x4 (A1:D3) block of spline coefficients; it becomes a
# 2x2 (A1:B2) block with new values, so remove columns C:D and row 3,
# then set the new values for the remaining A1:B2 cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused columns C and D, and row 3, shrinking the sheet
# from A1:D3 down to A1:B2.
$ws.Range("C1:D3").EntireColumn.Delete()
$ws.Range("A3:B3").EntireRow.Delete()

# Update the remaining 2x2 block of values.
$ws.Range("A1").Value = 3
$ws.Range("B1").Value = -2
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = -6
